# practica_4.xlsx - "se subio las nuevas clases"
#
# Changes applied:
#  1. Fix a typo in one of the instruction headers (shared string swap):
#     H6 -> "...y si encuentran ambas..." (now says "de la izquierda" instead of "a la izquierda")
#     K6 -> "...y si encuentran al menos una..." (unchanged wording)
#  2. Move the selection / viewport (sheetView) from C4/F7 to A5/A12.
#  3. Add the lookup formulas in columns F, H and K for the student rows (7-21).
#  4. Fix a data-entry typo in L12 (1029040136 -> 1029040139).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Header text fix (H6 / K6) -----------------------------------------
$ws.Range("H6").Value = "Busqueda de la cedulas de la izquierda,  y si encuentran ambas, mostrar busqueda finalizada, si no, decir siga intentando."
$ws.Range("K6").Value = "Busqueda de la cedulas a la izquierda,  y si encuentran al menos una, mostrar punto encontrado, si no encuentra ninguna, decir objetivo perdido."

# --- 2. Selection / view ----------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A12").Select()

# --- 3. Lookup formulas -----------------------------------------------------
$dRange = "`$D`$7,`$D`$8,`$D`$9,`$D`$10,`$D`$11,`$D`$12,`$D`$13,`$D`$14,`$D`$15,`$D`$16,`$D`$17,`$D`$18,`$D`$19,`$D`$20,`$D`$21,`$D`$22,`$D`$23,`$D`$24,`$D`$25,`$D`$26,`$D`$27,`$D`$28,`$D`$29,`$D`$30,`$D`$31,`$D`$32,`$D`$33,`$D`$34,`$D`$35"

$fFormula = '=IF(IFERROR(SEARCH(G7,CONCATENATE(' + $dRange + ')),0)>0,E7,"ninguno")'
$hFormula = '=IF(AND(IFERROR(SEARCH(I7,CONCATENATE(' + $dRange + ')),0)>0,IFERROR(SEARCH(J7,CONCATENATE(' + $dRange + ')),0)>0),"busqueda_finalizada","siga_intentando")'
$kFormula = '=IF(OR(IFERROR(SEARCH(L7,CONCATENATE(' + $dRange + ')),0)>0,IFERROR(SEARCH(M7,CONCATENATE(' + $dRange + ')),0)>0,IFERROR(SEARCH(N7,CONCATENATE(' + $dRange + ')),0)>0),"punto encontrado","objeto perdido")'

# Columns F and H only go down to row 13, column K goes down to row 21.
$ws.Range("F7:F13").Formula = $fFormula
$ws.Range("H7:H13").Formula = $hFormula
$ws.Range("K7:K21").Formula = $kFormula

# --- 4. Data fix -------------------------------------------------------------
$ws.Range("L12").Value = 1029040139
